# Apply updated N-gram statistics values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 387275

$ws.Range("D4").Value = 1.654744389679203
$ws.Range("E4").Value = 1.648895487702537

$ws.Range("D6").Value = 1.654744389679203
$ws.Range("E6").Value = 1.648895487702537

$ws.Range("D9").Value = 4.67844810573635
$ws.Range("E9").Value = 4.632017571323674
